$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"
$ws.Range("F19").Value = "-"

$ws.Range("B20").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"

$ws.Range("B21").Value = "[-, -, -, 'ELM-2NA-Sistemas de Refrigeração']"
$ws.Range("D21").Value = "-"
$ws.Range("F21").Value = "-"
